$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 17291.428
$ws.Range("I6").Value = 2810.25
$ws.Range("J6").Value = 36599.668
$ws.Range("K6").Value = 8430.75
$ws.Range("L6").Value = 109799.004
$ws.Range("M6").Value = -8318.75
$ws.Range("N6").Value = -110023.004
$ws.Range("H31").Value = 320
$ws.Range("I31").Value = 320
$ws.Range("K31").Value = 960
$ws.Range("M31").Value = -730
$ws.Range("H86").Value = 3954.65
$ws.Range("I86").Value = 3974.25
$ws.Range("K86").Value = 3974.25
$ws.Range("M86").Value = -2851.25
$ws.Range("H89").Value = 3954.65
$ws.Range("I89").Value = 3974.25
$ws.Range("K89").Value = 19871.25
$ws.Range("M89").Value = -14255.25
$ws.Range("H137").Value = 1572975.5
$ws.Range("I137").Value = 2404694.2
$ws.Range("J137").Value = 7387.353
$ws.Range("K137").Value = 7214082.600000001
$ws.Range("L137").Value = 22162.059
$ws.Range("M137").Value = -7211532.600000001
$ws.Range("N137").Value = -27262.059
$ws.Range("H138").Value = 1493.88
$ws.Range("I138").Value = 711.4167
$ws.Range("J138").Value = 1934.0156
$ws.Range("K138").Value = 2134.2501
$ws.Range("L138").Value = 5802.0468
$ws.Range("M138").Value = 3005.7499
$ws.Range("N138").Value = -16082.0468
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1997.2
$ws.Range("I61").Value = 1371.9166
$ws.Range("J61").Value = 4498.3335
$ws.Range("K61").Value = 1371.9166
$ws.Range("L61").Value = 4498.3335
$ws.Range("M61").Value = -1159.9166
$ws.Range("N61").Value = -4922.3335
$ws.Range("H74").Value = 1607.1282
$ws.Range("I74").Value = 1564.6538
$ws.Range("J74").Value = 1692.0769
$ws.Range("K74").Value = 1564.6538
$ws.Range("L74").Value = 1692.0769
$ws.Range("M74").Value = -690.6538
$ws.Range("N74").Value = -3440.0769
$ws.Range("H77").Value = 1607.1282
$ws.Range("I77").Value = 1564.6538
$ws.Range("J77").Value = 1692.0769
$ws.Range("K77").Value = 7823.269
$ws.Range("L77").Value = 8460.3845
$ws.Range("M77").Value = -3455.269
$ws.Range("N77").Value = -17196.3845
$ws.Range("H132").Value = 23811968
$ws.Range("I132").Value = 33334832
$ws.Range("J132").Value = 4810.1665
$ws.Range("K132").Value = 100004496
$ws.Range("L132").Value = 14430.4995
$ws.Range("M132").Value = -100001966
$ws.Range("N132").Value = -19490.4995
$ws.Range("H136").Value = 1997.2
$ws.Range("I136").Value = 1371.9166
$ws.Range("J136").Value = 4498.3335
$ws.Range("K136").Value = 4115.7498
$ws.Range("L136").Value = 13495.0005
$ws.Range("M136").Value = -1565.7498
$ws.Range("N136").Value = -18595.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2924.2788
$ws.Range("I134").Value = 1785.174
$ws.Range("J134").Value = 3613.7368
$ws.Range("K134").Value = 5355.522
$ws.Range("L134").Value = 10841.2104
$ws.Range("M134").Value = -2820.522
$ws.Range("N134").Value = -15911.2104
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1585.7567
$ws.Range("I58").Value = 1302.6
$ws.Range("J58").Value = 2799.2856
$ws.Range("K58").Value = 1302.6
$ws.Range("L58").Value = 2799.2856
$ws.Range("M58").Value = -1099.6
$ws.Range("N58").Value = -3205.2856
$ws.Range("H80").Value = 26654.375
$ws.Range("J80").Value = 26654.375
$ws.Range("L80").Value = 26654.375
$ws.Range("N80").Value = -28900.375
$ws.Range("H83").Value = 26654.375
$ws.Range("J83").Value = 26654.375
$ws.Range("L83").Value = 79963.125
$ws.Range("N83").Value = -91195.125
$ws.Range("H123").Value = 21473.684
$ws.Range("J123").Value = 21473.684
$ws.Range("L123").Value = 21473.684
$ws.Range("N123").Value = -31273.684
$ws.Range("H132").Value = 96007.2
$ws.Range("I132").Value = 1801.625
$ws.Range("J132").Value = 203670.72
$ws.Range("K132").Value = 5404.875
$ws.Range("L132").Value = 611012.16
$ws.Range("M132").Value = -2874.875
$ws.Range("N132").Value = -616072.16
$ws.Range("H136").Value = 1585.7567
$ws.Range("I136").Value = 1302.6
$ws.Range("J136").Value = 2799.2856
$ws.Range("K136").Value = 3907.8
$ws.Range("L136").Value = 8397.856800000001
$ws.Range("M136").Value = -1357.8
$ws.Range("N136").Value = -13497.8568
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 162326.2
$ws.Range("I56").Value = 162326.2
$ws.Range("K56").Value = 162326.2
$ws.Range("M56").Value = -161796.2
$ws.Range("H122").Value = 9629.166999999999
$ws.Range("J122").Value = 22551.8
$ws.Range("L122").Value = 202966.2
$ws.Range("N122").Value = -207866.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 37390
$ws.Range("J128").Value = 37390
$ws.Range("L128").Value = 37390
$ws.Range("N128").Value = -47350
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2788.2
$ws.Range("I132").Value = 2154.8147
$ws.Range("J132").Value = 3738.2778
$ws.Range("K132").Value = 6464.4441
$ws.Range("L132").Value = 11214.8334
$ws.Range("M132").Value = -3934.4441
$ws.Range("N132").Value = -16274.8334
$ws.Range("H136").Value = 2042.6207
$ws.Range("I136").Value = 1574.5454
$ws.Range("J136").Value = 3513.7144
$ws.Range("K136").Value = 4723.6362
$ws.Range("L136").Value = 10541.1432
$ws.Range("M136").Value = -2173.6362
$ws.Range("N136").Value = -15641.1432
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 10000758
$ws.Range("I107").Value = 528.5714
$ws.Range("J107").Value = 33334628
$ws.Range("K107").Value = 1585.7142
$ws.Range("L107").Value = 100003884
$ws.Range("M107").Value = 334.2857999999999
$ws.Range("N107").Value = -100007724
$ws.Range("H122").Value = 3175538
$ws.Range("I122").Value = 4762806
$ws.Range("K122").Value = 14288418
$ws.Range("M122").Value = -14285968
$ws.Range("H132").Value = 1892067.8
$ws.Range("I132").Value = 2289597.8
$ws.Range("J132").Value = 3799.5
$ws.Range("K132").Value = 6868793.399999999
$ws.Range("L132").Value = 11398.5
$ws.Range("M132").Value = -6866263.399999999
$ws.Range("N132").Value = -16458.5
$ws.Range("H136").Value = 1112303.9
$ws.Range("I136").Value = 1795760.2
$ws.Range("J136").Value = 1687.125
$ws.Range("K136").Value = 5387280.6
$ws.Range("L136").Value = 5061.375
$ws.Range("M136").Value = -5384730.6
$ws.Range("N136").Value = -10161.375
